$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.911.09'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.814.58'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.81'
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3657'
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07363'
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8700'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("D12").Value = '1.801.59'
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.363'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07092'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.504'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.35'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008700'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").Value = '26.952.30'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.61'
$ws.Range("E23").Value = '  -0.27%  '
$ws.Range("D24").Value = '2.020.76'
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.895'
$ws.Range("E25").Value = '  +0.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.80'
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.34'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.122'
$ws.Range("E28").Value = '  -1.47%  '
$ws.Range("E29").Value = '  -0.36%  '
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08882'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7554'
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.158'
$ws.Range("E33").Value = '  +0.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.480'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05287'
$ws.Range("E38").Value = '  +0.57%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.984'
$ws.Range("E40").Value = '  +1.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.253'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5305'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.326'
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1652'
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4869'
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.36'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.23'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.660'
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06296'
$ws.Range("E51").Value = '  +0.05%  '
